# "Fix new output table for Visual modeling and matrix"
#
# A new step ("Update predictive Model") needs to be listed in the "Steps"
# table, right above the existing "Update Workspace Plane Status" row, and
# the "Steps" sheet becomes the active/selected sheet of the workbook.

$wb = $excel.ActiveWorkbook
$wsSteps = $wb.Worksheets.Item("Steps")

# Insert a new row above the current row 6 ("Update Workspace Plane Status"),
# shifting the existing rows 6-9 down to 7-10.
$wsSteps.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the "Update predictive Model" step.
$wsSteps.Range("A6").Value = "Update_Nightly"
$wsSteps.Range("B6").Value = "Update predictive Model"
$wsSteps.Range("C6").Value = "Task"
$wsSteps.Range("D6").Value = "[Commons] Evaluate Formula"

# Match the (unstyled) formatting used by the rest of the table's data rows,
# rather than the formatting inherited from the row above on insert.
$wsSteps.Range("A6:D6").Style = $wsSteps.Range("A5:D5").Style

# The "Steps" sheet becomes the active sheet/tab, with cell B14 selected
# (and the previously active "Parameters" sheet loses that status).
$wsSteps.Activate()
$wsSteps.Range("B14").Select()
